$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Package info" sheet: bump version / timestamps for the v1.1.1 release
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Package info")
$wsInfo.Range("B4").Value = "V1.1.1"
$wsInfo.Range("B6").Value = "20210310T014404"
$wsInfo.Range("B7").Value = "COVAC_TRACKER_V1.1.1_DHIS2.34.3-80c86cc_20210310T014404"

# ---------------------------------------------------------------------------
# 2) "trackedEntityAttributes": rename the "COVAC - Sex" attribute code
# ---------------------------------------------------------------------------
$wsTea = $wb.Worksheets.Item("trackedEntityAttributes")
$wsTea.Range("B11").Value = "patinfo_sex"

# ---------------------------------------------------------------------------
# 3) "programs": bump the "Last updated" date for the tracked entity type
# ---------------------------------------------------------------------------
$wsPrograms = $wb.Worksheets.Item("programs")
$wsPrograms.Range("C2").Value = "2021-03-08"

# ---------------------------------------------------------------------------
# 4) "programRules": add a new program rule row, pushing the alphabetically
#    later rules down by one row
# ---------------------------------------------------------------------------
$wsRules = $wb.Worksheets.Item("programRules")
$wsRules.Rows.Item(21).Insert()
$wsRules.Cells.Item(21,1).Value = "R1bzqObecyQ"
$wsRules.Cells.Item(21,2).Value = "Hide Suggested date for next dose if vaccine product has no more doses"
$wsRules.Cells.Item(21,3).Value = "All vaccine types with two doses, after they receive one does, the ""next dose date"" will be hidden."
$wsRules.Cells.Item(21,5).Value = "yDuAzyqYABS"

# Restore the banded row shading (odd rows use the row-3 style, even rows
# the row-2 style) which Excel's row-insert would otherwise shift out of
# phase with the row index for every row from the insertion point on.
$wsRules.Range("A2:E2").Copy()
for ($r = 22; $r -le 56; $r += 2) {
    $wsRules.Range("A" + $r + ":E" + $r).PasteSpecial(-4122)
}
$wsRules.Range("A3:E3").Copy()
for ($r = 21; $r -le 55; $r += 2) {
    $wsRules.Range("A" + $r + ":E" + $r).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 5) "dataElements": give "COVAC - AEFIs present" a namespaced code
# ---------------------------------------------------------------------------
$wsDE = $wb.Worksheets.Item("dataElements")
$wsDE.Range("C2").Value = "COVAC_AEFIs_present"

# ---------------------------------------------------------------------------
# 6) "dataElementGroups": re-order the data elements belonging to the
#    "COVAC - Covid-19 vaccination registry" group (rows 17-21)
# ---------------------------------------------------------------------------
$wsDEG = $wb.Worksheets.Item("dataElementGroups")
$wsDEG.Range("B17").Value = "COVAC - Malignancy"
$wsDEG.Range("B18").Value = "COVAC - Dose Number"
$wsDEG.Range("B19").Value = "COVAC - Cardiovascular Disease"
$wsDEG.Range("B20").Value = "COVAC - Chronic Lung Disease"
$wsDEG.Range("B21").Value = "COVAC - AEFIs present"

# ---------------------------------------------------------------------------
# 7) "optionSets" / "options": correct the 3rd trimester week range
# ---------------------------------------------------------------------------
$wsOptionSets = $wb.Worksheets.Item("optionSets")
$wsOptionSets.Range("D5").Value = "1st Trimester (1-12 weeks); 2nd Trimester (13-28 weeks); 3rd Trimester (29-40 weeks)"

$wsOptions = $wb.Worksheets.Item("options")
$wsOptions.Range("B7").Value = "3rd Trimester (29-40 weeks)"

# ---------------------------------------------------------------------------
# 8) "visualizations": re-order the rows (row 6, "Number of doses
#    administered", keeps its place; all the others move around it)
# ---------------------------------------------------------------------------
$wsViz = $wb.Worksheets.Item("visualizations")

$wsViz.Range("A2").Value = "COVAC - People with completed vaccination schedule"
$wsViz.Range("B2").Value = "COVAC - People with completed vaccination schedule"
$wsViz.Range("D2").Value = "TWG0cq8P539"

$wsViz.Range("A3").Value = "COVAC - People receiving COV-2 vs People completing the vaccination schedule (Cov-C)"
$wsViz.Range("B3").Value = ""
$wsViz.Range("D3").Value = "wHd33PaphEC"

$wsViz.Range("A4").Value = "COVAC - At least one underlying condition"
$wsViz.Range("B4").Value = ""
$wsViz.Range("D4").Value = "gNsB9zivLTy"

$wsViz.Range("A5").Value = "COVAC - Underlying conditions"
$wsViz.Range("B5").Value = ""
$wsViz.Range("D5").Value = "vFkbMQiABfj"

$wsViz.Range("A7").Value = "COVAC - Complete vaccination uptake"
$wsViz.Range("B7").Value = ""
$wsViz.Range("D7").Value = "aUjo2Myd25f"

$wsViz.Range("A8").Value = "COVAC - Vaccine uptake by sex"
$wsViz.Range("B8").Value = ""
$wsViz.Range("D8").Value = "KV7fffdXnlY"

$wsViz.Range("A9").Value = "COVAC - Vaccine uptake by age group"
$wsViz.Range("B9").Value = ""
$wsViz.Range("D9").Value = "BWlYGFBDbO2"

$wsViz.Range("A10").Value = "COVAC - Vaccine uptake, last month"
$wsViz.Range("B10").Value = ""
$wsViz.Range("D10").Value = "MzSAvoJ0vLr"

$wsViz.Range("A11").Value = "COVAC - Drop-out from COV-1 to COV-c"
$wsViz.Range("B11").Value = ""
$wsViz.Range("D11").Value = "Hbs3xGj7XoN"

$wsViz.Range("A12").Value = "COVAC - Cumulative number of given doses, Last 4 weeks"
$wsViz.Range("B12").Value = "TEST TEST"
$wsViz.Range("D12").Value = "vmNUVdhuxN7"
